$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("worksheet")

# 1. Rename header B1 from "geo" to "geo_class"
$ws.Range("B1").Value = "geo_class"

# 2. Reorder rows 187:217 ("eth" block followed by "multi0" block) so that
#    the "multi0" rows (previously 201-217) come first, followed by the
#    "eth" rows (previously 187-200). This is a pure content swap of two
#    contiguous blocks within the same 31-row range - no rows are
#    inserted or deleted.
$rows = @(
    @('Aramaic', 'multi0', 'Syria, Iraq', $null, 'Syria, Iraq', $null),
    @('Central African', 'multi0', 'Angola,Cameroon,Central African Republic,Chad,Democratic Republic of the Congo,Eq. Guinea,Gabon,Republic of the Congo,Sao Tome and Principe', $null, 'Angola,Cameroon,Central African Republic,Chad,Democratic Republic of the Congo,Eq. Guinea,Gabon,Republic of the Congo,Sao Tome and Principe', $null),
    @('Congolese', 'multi0', 'Democratic Republic of the Congo, Republic of the Congo', $null, 'Democratic Republic of the Congo, Republic of the Congo', $null),
    @('East African', 'multi0', 'Burundi, Comoros, Djibouti, Eritrea, Ethiopia, Kenya, Madagascar, Malawi, Mauritius, Mozambique, Rwanda, S. Sudan, Seychelles, Somalia, Tanzania, Uganda, Zambia, Zimbabwe', $null, 'Burundi, Comoros, Djibouti, Eritrea, Ethiopia, Kenya, Madagascar, Malawi, Mauritius, Mozambique, Rwanda, S. Sudan, Seychelles, Somalia, Tanzania, Uganda, Zambia, Zimbabwe', $null),
    @('Finnish-Swedish', 'multi0', 'Finland, Sweden', $null, 'Finland, Sweden', $null),
    @('Korean', 'multi0', 'South Korea, North Korea', $null, 'South Korea, North Korea', $null),
    @('Melanesian', 'multi0', 'Papua New Guinea, Fiji, Vanuatu', $null, 'Papua New Guinea, Fiji, Vanuatu', $null),
    @('Micronesian', 'multi0', 'Micronesia, Marshall Islands, Palau', $null, 'Micronesia, Marshall Islands, Palau', $null),
    @('North African', 'multi0', 'Egypt, Libya, Tunisia, Algeria, Morocco', $null, 'Egypt, Libya, Tunisia, Algeria, Morocco', $null),
    @('North American Indian', 'multi0', 'United States, Canada', $null, 'United States, Canada', $null),
    @('Papuan', 'multi0', 'Papua New Guinea, Indonesia', $null, 'Papua New Guinea, Indonesia', $null),
    @('Persian Gulf', 'multi0', 'Iran, Iraq, Kuwait, Saudi Arabia, Bahrain, Qatar, United Arab Emirates, Oman', $null, 'Iran, Iraq, Kuwait, Saudi Arabia, Bahrain, Qatar, United Arab Emirates, Oman', $null),
    @('Polynesian', 'multi0', 'Samoa, Tonga, Tuvalu', $null, 'Samoa, Tonga, Tuvalu', $null),
    @('South American', 'multi0', 'Argentina, Bolivia, Brazil, Brazilian Island, Chile, Colombia, Ecuador, Falkland Is., Guyana, Paraguay, Peru, Suriname, Uruguay, Venezuela', $null, 'Argentina, Bolivia, Brazil, Brazilian Island, Chile, Colombia, Ecuador, Falkland Is., Guyana, Paraguay, Peru, Suriname, Uruguay, Venezuela', $null),
    @('Swahili', 'multi0', 'Kenya, Tanzania, Uganda, Mozambique', $null, 'Kenya, Tanzania, Uganda, Mozambique', $null),
    @('West African', 'multi0', 'Benin, Burkina Faso, Cabo Verde, Côte d''Ivoire, Gambia, Ghana, Guinea, Guinea-Bissau, Liberia, Mali, Mauritania, Niger, Nigeria, Saint Helena, Senegal, Sierra Leone, Togo', $null, 'Benin, Burkina Faso, Cabo Verde, Côte d''Ivoire, Gambia, Ghana, Guinea, Guinea-Bissau, Liberia, Mali, Mauritania, Niger, Nigeria, Saint Helena, Senegal, Sierra Leone, Togo', $null),
    @('West Indies', 'multi0', 'Anguilla, Antigua and Barb., Aruba, Bahamas, Barbados, British Virgin Is., Cayman Is., Cuba, Curaçao, Dominica, Dominican Republic, Grenada, Haiti, Jamaica, Montserrat, Puerto Rico, Saint Kitts and Nevis, Saint Lucia, Saint-Martin, Sint Maarten, St-Barthélemy, St. Vin. and Gren., Trinidad and Tobago, Turks and Caicos Is., U.S. Virgin Is.', $null, 'Anguilla, Antigua and Barb., Aruba, Bahamas, Barbados, British Virgin Is., Cayman Is., Cuba, Curaçao, Dominica, Dominican Republic, Grenada, Haiti, Jamaica, Montserrat, Puerto Rico, Saint Kitts and Nevis, Saint Lucia, Saint-Martin, Sint Maarten, St-Barthélemy, St. Vin. and Gren., Trinidad and Tobago, Turks and Caicos Is., U.S. Virgin Is.', $null),
    @('African American', 'eth', $null, $null, 'United States', 'African-American'),
    @('Altaic', 'eth', $null, $null, 'Russia, China, Mongolia', $null),
    @('Druze', 'eth', $null, $null, 'Lebanon, Syria, Israel', $null),
    @('Eskimo', 'eth', $null, $null, 'United States, Canada, Greenland, Russia', 'Inuit, Yupik'),
    @('Frisian', 'eth', $null, $null, 'Netherlands, Germany', $null),
    @('Gypsy', 'eth', $null, $null, 'Romania, Bulgaria, Hungary, Spain', $null),
    @('Kurdish', 'eth', $null, $null, 'Turkey, Iraq, Iran, Syria', 'Kurdistan'),
    @('Ladinian', 'eth', $null, $null, 'Italy', 'Ladin, Dolomites'),
    @('Lappish', 'eth', $null, $null, 'Norway, Sweden, Finland, Russia', $null),
    @('Mayan', 'eth', $null, $null, 'Mexico, Guatemala, Belize, Honduras, El Salvador', 'Mayan Empire'),
    @('Palestinian', 'eth', $null, $null, 'Palestine', $null),
    @('Sorbian', 'eth', $null, 'Slavic', 'Germany', 'Lusatia'),
    @('Spanish-American', 'eth', $null, $null, 'United States', $null),
    @('Tibetan', 'eth', $null, $null, 'China', 'Tibet')
)

$startRow = 187
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $vals = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
    $ws.Cells.Item($r, 5).Value = $vals[4]
    $ws.Cells.Item($r, 6).Value = $vals[5]
}

# 3. Restore the view's active cell/selection to where the user last
#    clicked after the reorder.
$ws.Activate() | Out-Null
$ws.Range("C204").Select() | Out-Null
